# Audit pass over the L-curve test workbook:
#  - remove the stray "Sheet" row (row 16) from optimization_parameters,
#    which also drops the now-unused "Sheet" shared string / number format.
#  - leave threshold_b as the active/selected sheet when done.

$wb = $excel.ActiveWorkbook

$wsOpt = $wb.Worksheets.Item("optimization_parameters")
$wsOpt.Activate() | Out-Null
$wsOpt.Rows.Item(16).Select() | Out-Null
$wsOpt.Rows.Item(16).Delete() | Out-Null

$wsThreshold = $wb.Worksheets.Item("threshold_b")
$wsThreshold.Activate() | Out-Null
